$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 2, shifting ALFONSO 50cl (and everything
# below it) down by one row. The new row 2 becomes a quantity-count row whose
# description is a literal single apostrophe character.
$ws.Rows(2).Insert()

# Populate the newly inserted row 2.
# Using "''" makes Excel's parser treat the first quote as the text-prefix
# marker and store a literal single apostrophe as the cell's text content.
$ws.Range("A2").Value = "''"
$ws.Range("B2").Value = 72
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 12
# Drop the quote-prefix cell style picked up from typing a leading apostrophe
# so the row matches the workbook's unstyled data rows.
$ws.Range("A2:D2").ClearFormats()

# Update the Case/Piece/Piece Per Case values that changed on the
# now-shifted-down rows (row numbers below are POST-insert row numbers).
$ws.Range("B3").Value = 23
$ws.Range("C3").Value = 11

$ws.Range("B4").Value = 36

$ws.Range("B5").Value = 32

$ws.Range("B6").Value = 1

$ws.Range("B11").Value = 317

$ws.Range("B17").Value = 197

$ws.Range("B18").Value = 29

$ws.Range("B21").Value = 24

$ws.Range("B22").Value = 0
$ws.Range("D22").Value = 12
